$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: add password value of literal `""`
$ws.Range("C5").Value = '""'

# Row 6: add username value of literal `""`
$ws.Range("B6").Value = '""'

# Row 7: add username and password values of literal `""`
$ws.Range("B7").Value = '""'
$ws.Range("C7").Value = '""'

# Rows 9-11: clear out the stray "NULL" username/password test entries
$ws.Range("C9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()

# Update the active selection to C5
$ws.Range("C5").Select()
